$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename project labels in column A
$ws.Range("A5").Value = "• Memory Bloo"
$ws.Range("A7").Value = "• Resonant Bloom"

# Update selection to reflect the active cell now being A7
$ws.Range("A7").Select()
